$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.285.90"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "3.818.07"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = $ws.Range("C4").Style
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'708.55"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").Value = "'172.03"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "3.817.96"
$ws.Range("E7").Value = "  -0.90%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'7.63"
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = "  +5.65%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "'36.07"
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "4.459.52"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "3.785.50"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "71.209.62"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "'17.53"
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'516.90"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = "  +3.67%  "
$ws.Range("D22").Value = "'10.71"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").Value = "'84.71"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").Value = "3.967.45"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'2.05"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("D31").Value = "'3.05"
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").Value = "'7.41"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("D34").Value = "'29.22"
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "3.786.14"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "'2.38"
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").Value = "'5.98"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("E42").Value = "  -2.48%  "
$ws.Range("D43").Value = "'3.29"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "'167.24"
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").Value = "'49.31"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").Value = "'422.74"
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("D50").Value = "'8.66"
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  +4.04%  "
